$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing last data row (244) as the style template for column A
# (date column, style index 2: bordered, bold, centered, custom date format).
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 245
$ws.Range("A245").Value = 44319
$ws.Range("B245").Value = 0
$ws.Range("C245").Value = 20
$ws.Range("D245").Value = 157.0475068708284

# Row 246
$ws.Range("A246").Value = 44320
$ws.Range("B246").Value = 0
$ws.Range("C246").Value = 18
$ws.Range("D246").Value = 141.3427561837456

# Row 247
$ws.Range("A247").Value = 44321
$ws.Range("B247").Value = 0
$ws.Range("C247").Value = 14
$ws.Range("D247").Value = 109.9332548095799
